$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was generically named "1" -- give it the proper municipality name.
$ws.Name = "Tkibuli"

# "upgrade left table until javakheti": the Urban/Rural rows for this
# municipality no longer carry any published per-year figures except a
# single confidential marker in the first (2010) column; the Total row
# keeps its three known data points (2010, 2012, 2018).
$ws.Range("B6").Value = "..."
$ws.Range("C6:O6").Value = "…"

$ws.Range("B7").Value = "..."
$ws.Range("C7:O7").Value = "…"

# Total row: only B5/D5/J5 stay numeric (unchanged); every other year column
# is unavailable.
$ws.Range("C5").Value = "…"
$ws.Range("E5:I5").Value = "…"
$ws.Range("K5:O5").Value = "…"

# The previously blank row 8 is removed, pulling the footnote up from row 9.
$ws.Rows.Item(8).Delete()
